$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.675.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.113.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5267"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4511"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.75"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09006"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.107.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.818"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.021"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001179"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.014"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.35"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.304"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.720.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.84%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356.03"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.534"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.80"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.346"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.021"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.899"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.19"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02633"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2308"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.58"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6867"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.279"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.320"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6417"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.764"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.244"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07283"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.49%  "

# Row 51: coin changed from WEMIXTOKEN to Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.05%  "
